# Weekly update for "Hortaliza, Vega Modelo de Temuco - Choclo": insert four
# new price rows (new reporting date 44578) before the current row 319, which
# pushes all the existing rows below it down by four (319-339 -> 323-343).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 319 (shifts old rows 319:339 -> 323:343)
$ws.Range("A319:A322").EntireRow.Insert()

# Fill the 4 newly inserted rows with this week's data. Columns A, B, C, E, F,
# G and R are constant for every row in this market/category block.
$newRows = @(
    @{ Row=319; D=44578; H="Choclero";           I="Primera"; J=50000; K=300; L=300; M=300; N="`$/unidad"; O="Región del Maule"; P=300; Q=1 },
    @{ Row=320; D=44578; H="Dulce o Americano";  I="Primera"; J=70000; K=200; L=200; M=200; N="`$/unidad"; O="Región del Maule"; P=200; Q=1 },
    @{ Row=321; D=44578; H="Dulce o Americano";  I="Segunda"; J=8000;  K=160; L=160; M=160; N="`$/unidad"; O="Región del Maule"; P=160; Q=1 },
    @{ Row=322; D=44578; H="Dulce o Americano";  I="Tercera"; J=500;   K=130; L=130; M=130; N="`$/unidad"; O="Región del Maule"; P=130; Q=1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 10
    $ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($row, 3).Value = "La Araucanía"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 9
    $ws.Cells.Item($row, 6).Value = 100112024
    $ws.Cells.Item($row, 7).Value = "Choclo"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
